#
# Small update to slides text
#
# 1. Refresh the cached "datetimeFigureOut" footer field text (25.09.2018 -> 19.11.2018)
#    on the slide master and every slide layout.
# 2. Slide 11 ("ASP.NET Session isn't a replacement for cache..."): extend the
#    bullet about in-memory caches with additional wording about cache invalidation.
# 3. Slide 14 ("Passing data among requests"): fix "follwoing" spacing and add
#    " in general" before "for more sophisticated applications".
# 4. Slide 4 ("StateServer"): clarify that StateServer uses a dedicated Win
#    Service *process* memory.
#

# A paragraph TextRange's .Text ends with a trailing CR (chr 13) for every
# paragraph that isn't the very last one in the text frame; that CR is not
# part of the addressable .Characters() range, so strip it before computing
# offsets/lengths.
function Get-VisibleText($range) {
    $txt = $range.Text
    if ($txt.Length -gt 0 -and [int][char]$txt[$txt.Length - 1] -eq 13) {
        return $txt.Substring(0, $txt.Length - 1)
    }
    return $txt
}

# Replace the (unique) substring $oldText inside $range with $newText,
# touching only the characters that actually changed so that unrelated runs
# in the same paragraph keep their original formatting.
function Replace-InRange($range, $oldText, $newText) {
    $visibleText = Get-VisibleText($range)
    $idx = $visibleText.IndexOf($oldText)
    if ($idx -lt 0) {
        throw ("Replace-InRange: substring not found: [" + $oldText + "] in [" + $visibleText + "]")
    }
    $startPos = $idx + 1
    $c = $range.Characters($startPos, $oldText.Length)
    $c.Text = $newText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date field text on the slide master + every layout.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "25.09.2018") {
                    $tr.Characters(1, $tr.Length).Text = "19.11.2018"
                }
            }
        }
    }
}

Update-DatePlaceholder($p.SlideMaster)

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li))
}

# ---------------------------------------------------------------------------
# 2) Slide 11 - extend the in-memory caches bullet.
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$tr11 = $slide11.Shapes.Item(1).TextFrame.TextRange
$para11_5 = $tr11.Paragraphs(5, 1)
Replace-InRange $para11_5 "duplication" "duplication and possible issues with cache invalidation."

# ---------------------------------------------------------------------------
# 3) Slide 14 - fix "follwoing options:" spacing and expand "in general".
# ---------------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$tr14 = $slide14.Shapes.Item(1).TextFrame.TextRange

$para14_2 = $tr14.Paragraphs(2, 1)
Replace-InRange $para14_2 "follwoing " "follwoing"
Replace-InRange $para14_2 "options:" " options:"

$para14_7 = $tr14.Paragraphs(7, 1)
Replace-InRange $para14_7 " in more sophisticated applications. Find more about Cross-Page Posting at " " in general for more sophisticated applications. Find more about Cross-Page Posting at "

# ---------------------------------------------------------------------------
# 4) Slide 4 - StateServer description wording tweak.
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tr4 = $slide4.Shapes.Item(1).TextFrame.TextRange
$para4_1 = $tr4.Paragraphs(1, 1)
$oldRun4 = " " + [char]0x2013 + " saves session data in dedicated Win Service memory"
$newRun4 = " " + [char]0x2013 + " saves session data in a dedicated Win Service process memory"
Replace-InRange $para4_1 $oldRun4 $newRun4
